# Forming the consolidated report: fill in the "Absent" (column H) values
# that were previously left blank / incorrect for this student's rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows where H changes from 0 -> 1
$ws.Range("H3").Value  = 1
$ws.Range("H9").Value  = 1
$ws.Range("H12").Value = 1

# Rows where H changes from an empty inline string -> numeric 0
$ws.Range("H6").Value  = 0
$ws.Range("H10").Value = 0
$ws.Range("H14").Value = 0
